$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume figures (and two ranking swaps)
# as produced by the scheduled GitHub Actions refresh job.
# Cells whose new text would otherwise be auto-detected by Excel as a
# number (e.g. "1.00", "314.74") are written with the cell pre-formatted
# as Text so the literal string is preserved, matching the source feed.

$ws.Cells.Item(2, 4).Value = '44.610.39'
$ws.Cells.Item(2, 5).Value = '  +3.89%  '
$ws.Cells.Item(3, 4).Value = '2.422.26'
$ws.Cells.Item(3, 5).Value = '  +2.60%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '314.74'
$ws.Cells.Item(5, 5).Value = '  +3.88%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '101.07'
$ws.Cells.Item(6, 5).Value = '  +5.82%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.513'
$ws.Cells.Item(7, 5).Value = '  +2.52%  '
$ws.Cells.Item(8, 5).Value = '  -0.08%  '
$ws.Cells.Item(9, 5).Value = '  +8.65%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '35.42'
$ws.Cells.Item(10, 5).Value = '  +4.03%  '
$ws.Cells.Item(11, 5).Value = '  +2.03%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '18.98'
$ws.Cells.Item(12, 5).Value = '  +2.91%  '
$ws.Cells.Item(13, 5).Value = '  -2.16%  '
$ws.Cells.Item(14, 5).Value = '  +3.88%  '
$ws.Cells.Item(15, 4).Value = '2.800.71'
$ws.Cells.Item(15, 5).Value = '  +2.70%  '
$ws.Cells.Item(16, 4).Value = '2.395.77'
$ws.Cells.Item(16, 5).Value = '  +2.01%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.832'
$ws.Cells.Item(17, 5).Value = '  +5.27%  '
$ws.Cells.Item(18, 4).Value = '44.497.72'
$ws.Cells.Item(18, 5).Value = '  +3.68%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '12.48'
$ws.Cells.Item(19, 5).Value = '  +4.85%  '
$ws.Cells.Item(20, 5).Value = '  +2.37%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0922'
$ws.Cells.Item(21, 5).Value = '  +4.37%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '68.72'
$ws.Cells.Item(22, 5).Value = '  +1.12%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '242.76'
$ws.Cells.Item(23, 5).Value = '  +3.33%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.28'
$ws.Cells.Item(24, 5).Value = '  +5.42%  '
$ws.Cells.Item(25, 5).Value = '  +1.97%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.999'
$ws.Cells.Item(26, 5).Value = '  -0.14%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '25.24'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.25'
$ws.Cells.Item(28, 5).Value = '  -5.02%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '9.54'
$ws.Cells.Item(29, 5).Value = '  +2.55%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '33.26'
$ws.Cells.Item(30, 5).Value = '  +3.95%  '
$ws.Cells.Item(31, 5).Value = '  +1.44%  '
$ws.Cells.Item(32, 5).Value = '  +20.43%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '19.40'
$ws.Cells.Item(33, 5).Value = '  +10.85%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.18'
$ws.Cells.Item(34, 5).Value = '  +3.67%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0776'
$ws.Cells.Item(35, 5).Value = '  +8.54%  '
$ws.Cells.Item(36, 5).Value = '  +0.18%  '
$ws.Cells.Item(37, 5).Value = '  +3.10%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '4.50'
$ws.Cells.Item(38, 5).Value = '  +4.41%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.87'
$ws.Cells.Item(39, 5).Value = '  +1.82%  '
$ws.Cells.Item(40, 2).Value = 'Monero'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '119.90'
$ws.Cells.Item(40, 5).Value = '  -6.70%  '
$ws.Cells.Item(41, 5).Value = '  +1.84%  '
$ws.Cells.Item(42, 2).Value = 'WEMIXToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.20'
$ws.Cells.Item(42, 5).Value = '  -2.52%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '21.05'
$ws.Cells.Item(43, 5).Value = '  -0.96%  '
$ws.Cells.Item(44, 5).Value = '  +4.58%  '
$ws.Cells.Item(45, 4).Value = '1.941.03'
$ws.Cells.Item(45, 5).Value = '  +0.67%  '
$ws.Cells.Item(46, 5).Value = '  +1.60%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.95'
$ws.Cells.Item(47, 5).Value = '  +9.25%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '9.44'
$ws.Cells.Item(48, 5).Value = '  +2.34%  '
$ws.Cells.Item(49, 5).Value = '  +11.43%  '
$ws.Cells.Item(50, 2).Value = 'BitcoinSV'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '75.65'
$ws.Cells.Item(50, 5).Value = '  +5.92%  '
$ws.Cells.Item(51, 2).Value = 'MultiversX'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '54.48'
$ws.Cells.Item(51, 5).Value = '  +6.54%  '
